# Sassuolo report: "aggiornamento fino a 8/12" (update through 2021-12-08)
# Appends 79 new daily rows (386-464) to the existing data table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A carries a dedicated date-style format (thin border, bold, centered,
# custom date number format). Extend that formatting from the last existing
# row (A385) down through the newly appended rows before writing values, so the
# new cells reuse the same style instead of creating new ones.
$ws.Cells.Item(385, 1).Copy()
$ws.Range($ws.Cells.Item(386, 1), $ws.Cells.Item(464, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Bulk-write the new data (date serial, nuovi pos., somma mobile 7gg.,
# somma mobile 7gg. per 100mila abitanti) for rows 386-464 in one shot.
$data = New-Object 'object[,]' 79,4
$data[0,0] = 44460
$data[0,1] = 2
$data[0,2] = 18
$data[0,3] = 44.68940861015939
$data[1,0] = 44461
$data[1,1] = 0
$data[1,2] = 18
$data[1,3] = 44.68940861015939
$data[2,0] = 44462
$data[2,1] = 2
$data[2,2] = 11
$data[2,3] = 27.31019415065296
$data[3,0] = 44463
$data[3,1] = 0
$data[3,2] = 9
$data[3,3] = 22.34470430507969
$data[4,0] = 44464
$data[4,1] = 1
$data[4,2] = 8
$data[4,3] = 19.86195938229307
$data[5,0] = 44465
$data[5,1] = 0
$data[5,2] = 7
$data[5,3] = 17.37921445950643
$data[6,0] = 44466
$data[6,1] = 0
$data[6,2] = 5
$data[6,3] = 12.41372461393316
$data[7,0] = 44467
$data[7,1] = 1
$data[7,2] = 4
$data[7,3] = 9.930979691146533
$data[8,0] = 44468
$data[8,1] = 1
$data[8,2] = 5
$data[8,3] = 12.41372461393316
$data[9,0] = 44469
$data[9,1] = 1
$data[9,2] = 4
$data[9,3] = 9.930979691146533
$data[10,0] = 44470
$data[10,1] = 2
$data[10,2] = 6
$data[10,3] = 14.8964695367198
$data[11,0] = 44471
$data[11,1] = 1
$data[11,2] = 6
$data[11,3] = 14.8964695367198
$data[12,0] = 44472
$data[12,1] = 0
$data[12,2] = 6
$data[12,3] = 14.8964695367198
$data[13,0] = 44473
$data[13,1] = 2
$data[13,2] = 8
$data[13,3] = 19.86195938229307
$data[14,0] = 44474
$data[14,1] = 1
$data[14,2] = 8
$data[14,3] = 19.86195938229307
$data[15,0] = 44475
$data[15,1] = 1
$data[15,2] = 8
$data[15,3] = 19.86195938229307
$data[16,0] = 44476
$data[16,1] = 5
$data[16,2] = 12
$data[16,3] = 29.79293907343959
$data[17,0] = 44477
$data[17,1] = 1
$data[17,2] = 11
$data[17,3] = 27.31019415065296
$data[18,0] = 44478
$data[18,1] = 0
$data[18,2] = 10
$data[18,3] = 24.82744922786633
$data[19,0] = 44479
$data[19,1] = 2
$data[19,2] = 12
$data[19,3] = 29.79293907343959
$data[20,0] = 44480
$data[20,1] = 0
$data[20,2] = 10
$data[20,3] = 24.82744922786633
$data[21,0] = 44481
$data[21,1] = 2
$data[21,2] = 11
$data[21,3] = 27.31019415065296
$data[22,0] = 44482
$data[22,1] = 1
$data[22,2] = 11
$data[22,3] = 27.31019415065296
$data[23,0] = 44483
$data[23,1] = 1
$data[23,2] = 7
$data[23,3] = 17.37921445950643
$data[24,0] = 44484
$data[24,1] = 1
$data[24,2] = 7
$data[24,3] = 17.37921445950643
$data[25,0] = 44485
$data[25,1] = 0
$data[25,2] = 7
$data[25,3] = 17.37921445950643
$data[26,0] = 44486
$data[26,1] = 2
$data[26,2] = 7
$data[26,3] = 17.37921445950643
$data[27,0] = 44487
$data[27,1] = 0
$data[27,2] = 7
$data[27,3] = 17.37921445950643
$data[28,0] = 44488
$data[28,1] = 1
$data[28,2] = 6
$data[28,3] = 14.8964695367198
$data[29,0] = 44489
$data[29,1] = 0
$data[29,2] = 5
$data[29,3] = 12.41372461393316
$data[30,0] = 44490
$data[30,1] = 0
$data[30,2] = 4
$data[30,3] = 9.930979691146533
$data[31,0] = 44491
$data[31,1] = 2
$data[31,2] = 5
$data[31,3] = 12.41372461393316
$data[32,0] = 44492
$data[32,1] = 0
$data[32,2] = 5
$data[32,3] = 12.41372461393316
$data[33,0] = 44493
$data[33,1] = 0
$data[33,2] = 3
$data[33,3] = 7.448234768359899
$data[34,0] = 44494
$data[34,1] = 2
$data[34,2] = 5
$data[34,3] = 12.41372461393316
$data[35,0] = 44495
$data[35,1] = 6
$data[35,2] = 10
$data[35,3] = 24.82744922786633
$data[36,0] = 44496
$data[36,1] = 0
$data[36,2] = 10
$data[36,3] = 24.82744922786633
$data[37,0] = 44497
$data[37,1] = 2
$data[37,2] = 12
$data[37,3] = 29.79293907343959
$data[38,0] = 44498
$data[38,1] = 2
$data[38,2] = 12
$data[38,3] = 29.79293907343959
$data[39,0] = 44499
$data[39,1] = 2
$data[39,2] = 14
$data[39,3] = 34.75842891901286
$data[40,0] = 44500
$data[40,1] = 3
$data[40,2] = 17
$data[40,3] = 42.20666368737276
$data[41,0] = 44501
$data[41,1] = 3
$data[41,2] = 18
$data[41,3] = 44.68940861015939
$data[42,0] = 44502
$data[42,1] = 2
$data[42,2] = 14
$data[42,3] = 34.75842891901286
$data[43,0] = 44503
$data[43,1] = 0
$data[43,2] = 14
$data[43,3] = 34.75842891901286
$data[44,0] = 44504
$data[44,1] = 2
$data[44,2] = 14
$data[44,3] = 34.75842891901286
$data[45,0] = 44505
$data[45,1] = 1
$data[45,2] = 13
$data[45,3] = 32.27568399622623
$data[46,0] = 44506
$data[46,1] = 1
$data[46,2] = 12
$data[46,3] = 29.79293907343959
$data[47,0] = 44507
$data[47,1] = 2
$data[47,2] = 11
$data[47,3] = 27.31019415065296
$data[48,0] = 44508
$data[48,1] = 2
$data[48,2] = 10
$data[48,3] = 24.82744922786633
$data[49,0] = 44509
$data[49,1] = 1
$data[49,2] = 9
$data[49,3] = 22.34470430507969
$data[50,0] = 44510
$data[50,1] = 1
$data[50,2] = 10
$data[50,3] = 24.82744922786633
$data[51,0] = 44511
$data[51,1] = 6
$data[51,2] = 14
$data[51,3] = 34.75842891901286
$data[52,0] = 44512
$data[52,1] = 1
$data[52,2] = 14
$data[52,3] = 34.75842891901286
$data[53,0] = 44513
$data[53,1] = 6
$data[53,2] = 19
$data[53,3] = 47.17215353294603
$data[54,0] = 44514
$data[54,1] = 4
$data[54,2] = 21
$data[54,3] = 52.13764337851929
$data[55,0] = 44515
$data[55,1] = 4
$data[55,2] = 23
$data[55,3] = 57.10313322409256
$data[56,0] = 44516
$data[56,1] = 47
$data[56,2] = 69
$data[56,3] = 171.3093996722777
$data[57,0] = 44517
$data[57,1] = 1
$data[57,2] = 69
$data[57,3] = 171.3093996722777
$data[58,0] = 44518
$data[58,1] = 8
$data[58,2] = 71
$data[58,3] = 176.2748895178509
$data[59,0] = 44519
$data[59,1] = 9
$data[59,2] = 79
$data[59,3] = 196.136848900144
$data[60,0] = 44520
$data[60,1] = 9
$data[60,2] = 82
$data[60,3] = 203.5850836685039
$data[61,0] = 44521
$data[61,1] = 6
$data[61,2] = 84
$data[61,3] = 208.5505735140772
$data[62,0] = 44522
$data[62,1] = 5
$data[62,2] = 85
$data[62,3] = 211.0333184368638
$data[63,0] = 44523
$data[63,1] = 2
$data[63,2] = 40
$data[63,3] = 99.30979691146531
$data[64,0] = 44524
$data[64,1] = 22
$data[64,2] = 61
$data[64,3] = 151.4474402899846
$data[65,0] = 44525
$data[65,1] = 15
$data[65,2] = 68
$data[65,3] = 168.826654749491
$data[66,0] = 44526
$data[66,1] = 15
$data[66,2] = 74
$data[66,3] = 183.7231242862108
$data[67,0] = 44527
$data[67,1] = 2
$data[67,2] = 67
$data[67,3] = 166.3439098267044
$data[68,0] = 44528
$data[68,1] = 12
$data[68,2] = 73
$data[68,3] = 181.2403793634242
$data[69,0] = 44529
$data[69,1] = 5
$data[69,2] = 73
$data[69,3] = 181.2403793634242
$data[70,0] = 44530
$data[70,1] = 5
$data[70,2] = 76
$data[70,3] = 188.6886141317841
$data[71,0] = 44531
$data[71,1] = 0
$data[71,2] = 54
$data[71,3] = 134.0682258304782
$data[72,0] = 44532
$data[72,1] = 10
$data[72,2] = 49
$data[72,3] = 121.654501216545
$data[73,0] = 44533
$data[73,1] = 6
$data[73,2] = 40
$data[73,3] = 99.30979691146531
$data[74,0] = 44534
$data[74,1] = 11
$data[74,2] = 49
$data[74,3] = 121.654501216545
$data[75,0] = 44535
$data[75,1] = 8
$data[75,2] = 45
$data[75,3] = 111.7235215253985
$data[76,0] = 44536
$data[76,1] = 7
$data[76,2] = 47
$data[76,3] = 116.6890113709718
$data[77,0] = 44537
$data[77,1] = 8
$data[77,2] = 50
$data[77,3] = 124.1372461393317
$data[78,0] = 44538
$data[78,1] = 0
$data[78,2] = 50
$data[78,3] = 124.1372461393317

$ws.Range($ws.Cells.Item(386, 1), $ws.Cells.Item(464, 4)).Value = $data

